$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Age value for row 6 (James/Rachel entry) from 25 to 40
$ws.Range("C6").Value = 40

# Move the active selection from D6 to C6
$ws.Range("C6").Select()
